# Updated cryptos list price/volume figures.
# Column D ("Price") and column E ("Volume(1h)") are plain text cells in the
# source sheet (prices use "." as a thousands separator in several rows, e.g.
# "66.741.07", so the whole column is stored as text, not numbers). When a
# "Price" value looks like a bare decimal number (e.g. "580.06") we have to
# force the destination cell to Text format before writing it, otherwise
# Excel's COM layer auto-converts the string to a real number and mangles it
# (drops trailing zeros, "580.00" -> 580, "1.00" -> 1, etc). ClearFormats()
# afterwards drops the now-unneeded "@" number format again so the cell's
# style stays the same as before (the stored value itself, already text,
# is unaffected by ClearFormats).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$updates = @(
    @{ Row = 2;  Price = "66.749.18"; Volume = "+2.05%" },
    @{ Row = 3;  Price = "3.091.10";  Volume = "+5.28%" },
    @{ Row = 4;  Price = $null;       Volume = "+0.03%" },
    @{ Row = 5;  Price = "580.00";    Volume = "+2.02%" },
    @{ Row = 6;  Price = "168.34";    Volume = "+6.24%" },
    @{ Row = 7;  Price = "0.999";     Volume = "+0.03%" },
    @{ Row = 8;  Price = "3.087.33";  Volume = "+5.25%" },
    @{ Row = 9;  Price = $null;       Volume = "+1.69%" },
    @{ Row = 10; Price = "6.59";      Volume = "-2.19%" },
    @{ Row = 11; Price = $null;       Volume = "+3.64%" },
    @{ Row = 12; Price = "0.482";     Volume = "+4.93%" },
    @{ Row = 13; Price = $null;       Volume = "+2.56%" },
    @{ Row = 14; Price = "36.42";     Volume = "+6.11%" },
    @{ Row = 15; Price = $null;       Volume = "-0.53%" },
    @{ Row = 16; Price = "3.601.65";  Volume = "+5.23%" },
    @{ Row = 17; Price = "66.722.90"; Volume = "+2.07%" },
    @{ Row = 18; Price = "7.20";      Volume = "+3.07%" },
    @{ Row = 19; Price = "3.092.00";  Volume = "+5.45%" },
    @{ Row = 20; Price = "16.28";     Volume = "+3.83%" },
    @{ Row = 21; Price = "466.66";    Volume = "+5.09%" },
    @{ Row = 22; Price = $null;       Volume = "+3.41%" },
    @{ Row = 23; Price = $null;       Volume = "+3.39%" },
    @{ Row = 24; Price = "83.96";     Volume = "+2.13%" },
    @{ Row = 25; Price = $null;       Volume = "+5.40%" },
    @{ Row = 26; Price = "13.06";     Volume = "+7.82%" },
    @{ Row = 27; Price = "10.11";     Volume = "+0.64%" },
    @{ Row = 28; Price = $null;       Volume = "-0.04%" },
    @{ Row = 29; Price = "8.03";      Volume = "+0.04%" },
    @{ Row = 30; Price = $null;       Volume = "+2.14%" },
    @{ Row = 31; Price = $null;       Volume = "+3.99%" },
    @{ Row = 32; Price = $null;       Volume = "+1.16%" },
    @{ Row = 33; Price = "28.30";     Volume = "+4.69%" },
    @{ Row = 34; Price = $null;       Volume = "+3.75%" },
    @{ Row = 35; Price = $null;       Volume = "+0.10%" },
    @{ Row = 36; Price = "1.00";      Volume = "+3.24%" },
    @{ Row = 37; Price = "5.90";      Volume = "+2.75%" },
    @{ Row = 38; Price = "47.31";     Volume = "+5.37%" },
    @{ Row = 39; Price = $null;       Volume = "+6.34%" },
    @{ Row = 40; Price = $null;       Volume = "+6.15%" },
    @{ Row = 41; Price = "50.32";     Volume = "+1.37%" },
    @{ Row = 42; Price = $null;       Volume = "+1.64%" },
    @{ Row = 43; Price = $null;       Volume = "+2.48%" },
    @{ Row = 44; Price = $null;       Volume = "-0.36%" },
    @{ Row = 45; Price = $null;       Volume = "+2.79%" },
    @{ Row = 46; Price = "382.53";    Volume = "-0.18%" },
    @{ Row = 47; Price = "2.787.21";  Volume = "+3.32%" },
    @{ Row = 48; Price = "135.04";    Volume = "+1.20%" },
    @{ Row = 50; Price = "24.92";     Volume = "+6.79%" },
    @{ Row = 51; Price = $null;       Volume = "+1.77%" }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.Price) {
        $dRange = $ws.Range("D$row")
        # A decimal-looking price (no second "." thousands separator) would
        # otherwise be silently re-typed as a Number by Excel.
        if ($u.Price -match '^[0-9]+\.[0-9]+$') {
            Set-TextValue $dRange $u.Price
        } else {
            $dRange.Value = $u.Price
        }
    }

    if ($null -ne $u.Volume) {
        $ws.Range("E$row").Value = "  " + $u.Volume + "  "
    }
}
